$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-15
$values = @(
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(2, 2),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(3, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
